$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.798171611121401
$ws.Range("D2").Value = 4.474690696834815
$ws.Range("E2").Value = 16.48531796207418
$ws.Range("F2").Value = 26.05631150714657
$ws.Range("G2").Value = 32.62165228357534
$ws.Range("H2").Value = 14.62716152995666
$ws.Range("K2").Value = 11.80815382165913

$ws.Range("B3").Value = 7.726857158111801
$ws.Range("D3").Value = 4.488718641446136
$ws.Range("E3").Value = 15.548202864006
$ws.Range("F3").Value = 25.87975062817515
$ws.Range("G3").Value = 32.21707872641259
$ws.Range("H3").Value = 14.64181875524063
$ws.Range("K3").Value = 11.1729203446388

$ws.Range("B4").Value = 7.684635076878071
$ws.Range("D4").Value = 4.497648392147862
$ws.Range("E4").Value = 14.94846442266241
$ws.Range("F4").Value = 25.78087893220476
$ws.Range("G4").Value = 31.98154492478526
$ws.Range("H4").Value = 14.65484605124659
$ws.Range("K4").Value = 10.7617040526054

$ws.Range("B5").Value = 7.66784154203761
$ws.Range("D5").Value = 4.50136753528077
$ws.Range("E5").Value = 14.69821475571074
$ws.Range("F5").Value = 25.74301672522822
$ws.Range("G5").Value = 31.88891547326455
$ws.Range("H5").Value = 14.66116283631637
$ws.Range("K5").Value = 10.58887761902343

$ws.Range("B6").Value = 7.66507842184571
$ws.Range("D6").Value = 4.501989957681123
$ws.Range("E6").Value = 14.65631653858642
$ws.Range("F6").Value = 25.73687721789727
$ws.Range("G6").Value = 31.87374006130931
$ws.Range("H6").Value = 14.6622724695198
$ws.Range("K6").Value = 10.5598651667472

$ws.Range("B7").Value = 7.684406899847426
$ws.Range("D7").Value = 4.497698224341965
$ws.Range("E7").Value = 14.94511276330694
$ws.Range("F7").Value = 25.78035843882758
$ws.Range("G7").Value = 31.98028197358042
$ws.Range("H7").Value = 14.65492716678574
$ws.Range("K7").Value = 10.75939441241652

$ws.Range("B8").Value = 7.773268965766825
$ws.Range("D8").Value = 4.479462195153734
$ws.Range("E8").Value = 16.16740030896938
$ws.Range("F8").Value = 25.99347376072951
$ws.Range("G8").Value = 32.47955470377534
$ws.Range("H8").Value = 14.63137668487632
$ws.Range("K8").Value = 11.59356043256521

$ws.Range("B9").Value = 7.959090527140051
$ws.Range("D9").Value = 4.446185947940724
$ws.Range("E9").Value = 18.44633017329961
$ws.Range("F9").Value = 26.48547445317219
$ws.Range("G9").Value = 33.55529405824703
$ws.Range("H9").Value = 14.61734305045017
$ws.Range("K9").Value = 13.05884261341183

$ws.Range("B10").Value = 8.101475941776748
$ws.Range("D10").Value = 4.42321532765346
$ws.Range("E10").Value = 20.0955935661945
$ws.Range("F10").Value = 26.8896427146984
$ws.Range("G10").Value = 34.3966814822108
$ws.Range("H10").Value = 14.62686046681531
$ws.Range("K10").Value = 14.02902604570867

$ws.Range("B11").Value = 8.167263434401081
$ws.Range("D11").Value = 4.413078446028126
$ws.Range("E11").Value = 20.80406955884616
$ws.Range("F11").Value = 27.08216538860639
$ws.Range("G11").Value = 34.78874967107185
$ws.Range("H11").Value = 14.63553319862568
$ws.Range("K11").Value = 14.44698255130692

$ws.Range("B12").Value = 8.192299432714384
$ws.Range("D12").Value = 4.409284235438133
$ws.Range("E12").Value = 21.06638142871904
$ws.Range("F12").Value = 27.15625766793558
$ws.Range("G12").Value = 34.93840246873502
$ws.Range("H12").Value = 14.6394439523078
$ws.Range("K12").Value = 14.60187225344108

$ws.Range("B13").Value = 8.186902336320408
$ws.Range("D13").Value = 4.410099420047495
$ws.Range("E13").Value = 21.01015294263207
$ws.Range("F13").Value = 27.14024867976736
$ws.Range("G13").Value = 34.90612182182938
$ws.Range("H13").Value = 14.63857380845124
$ws.Range("K13").Value = 14.56866460517253

$ws.Range("B14").Value = 8.1693208097743
$ws.Range("D14").Value = 4.412765407123478
$ws.Range("E14").Value = 20.8257696109346
$ws.Range("F14").Value = 27.08823753501618
$ws.Range("G14").Value = 34.80103883373521
$ws.Range("H14").Value = 14.63584237167604
$ws.Range("K14").Value = 14.45979327328498

$ws.Range("B15").Value = 8.15856706395628
$ws.Range("D15").Value = 4.414404168380011
$ws.Range("E15").Value = 20.71205273977241
$ws.Range("F15").Value = 27.05653216854545
$ws.Range("G15").Value = 34.73682220187657
$ws.Range("H15").Value = 14.63425093657013
$ws.Range("K15").Value = 14.39266563590206

$ws.Range("B16").Value = 8.097195101632455
$ws.Range("D16").Value = 4.423884039900339
$ws.Range("E16").Value = 20.04845527816319
$ws.Range("F16").Value = 26.87723038380184
$ws.Range("G16").Value = 34.37123317938481
$ws.Range("H16").Value = 14.62638126735867
$ws.Range("K16").Value = 14.00123930828737

$ws.Range("B17").Value = 8.059789628217139
$ws.Range("D17").Value = 4.42977930185056
$ws.Range("E17").Value = 19.63068370883585
$ws.Range("F17").Value = 26.76941357895621
$ws.Range("G17").Value = 34.14923268294949
$ws.Range("H17").Value = 14.62266751833315
$ws.Range("K17").Value = 13.75510726423742

$ws.Range("B18").Value = 8.038371929024416
$ws.Range("D18").Value = 4.433199552609126
$ws.Range("E18").Value = 19.38645769528683
$ws.Range("F18").Value = 26.7082194217752
$ws.Range("G18").Value = 34.02243017468152
$ws.Range("H18").Value = 14.62094018890907
$ws.Range("K18").Value = 13.6113382148869

$ws.Range("B19").Value = 8.031137594376185
$ws.Range("D19").Value = 4.434362665807535
$ws.Range("E19").Value = 19.30309006027881
$ws.Range("F19").Value = 26.68764254424135
$ws.Range("G19").Value = 33.97965394590413
$ws.Range("H19").Value = 14.62042547112491
$ws.Range("K19").Value = 13.56228323310836

$ws.Range("B20").Value = 8.063761633728893
$ws.Range("D20").Value = 4.429148697016365
$ws.Range("E20").Value = 19.67556318860938
$ws.Range("F20").Value = 26.78080648560176
$ws.Range("G20").Value = 34.17277440395993
$ws.Range("H20").Value = 14.62302053578007
$ws.Range("K20").Value = 13.78153632851103

$ws.Range("B21").Value = 8.174481751016938
$ws.Range("D21").Value = 4.411981141159503
$ws.Range("E21").Value = 20.88008927704987
$ws.Range("F21").Value = 27.10348269919215
$ws.Range("G21").Value = 34.83187332695053
$ws.Range("H21").Value = 14.63662764256152
$ws.Range("K21").Value = 14.49186331052902

$ws.Range("B22").Value = 8.247554461078463
$ws.Range("D22").Value = 4.401019743816533
$ws.Range("E22").Value = 21.63252014740533
$ws.Range("F22").Value = 27.32126674289191
$ws.Range("G22").Value = 35.26946834846172
$ws.Range("H22").Value = 14.64917335966419
$ws.Range("K22").Value = 14.93638928895666

$ws.Range("B23").Value = 8.208496520258121
$ws.Range("D23").Value = 4.406846562001201
$ws.Range("E23").Value = 21.23410579946228
$ws.Range("F23").Value = 27.2044201120642
$ws.Range("G23").Value = 35.03534133870266
$ws.Range("H23").Value = 14.6421427411054
$ws.Range("K23").Value = 14.70094587966283

$ws.Range("B24").Value = 8.061965616712275
$ws.Range("D24").Value = 4.429433696862984
$ws.Range("E24").Value = 19.65528577385725
$ws.Range("F24").Value = 26.77565328369553
$ws.Range("G24").Value = 34.16212859927397
$ws.Range("H24").Value = 14.62285966657836
$ws.Range("K24").Value = 13.76959480047435

$ws.Range("B25").Value = 7.907703759762687
$ws.Range("D25").Value = 4.454925957908091
$ws.Range("E25").Value = 17.80168581028119
$ws.Range("F25").Value = 26.34468124832206
$ws.Range("G25").Value = 33.25472055204132
$ws.Range("H25").Value = 14.61767192313213
$ws.Range("K25").Value = 12.6809895545779

